# fix : correction gestion du fichier excel deja ouvert
# Refresh the crypto "Cours" (price) column and the dependent Profits /
# Profits % / Variation 24h columns, plus the last-update timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cryptos")

$newUpdate = "02/04/21 23:25"

# Row 4 - Bitcoin
$ws.Range("I4").Value = 50160.92
$ws.Range("L4").Value = -0.0046
$ws.Range("N4").Value = $newUpdate

# Row 5 - Etherum
$ws.Range("I5").Value = 1780.45
$ws.Range("L5").Value = 0.0535
$ws.Range("N5").Value = $newUpdate

# Row 6 - 0x (ZRX) : only Variation 24h / Update change
$ws.Range("L6").Value = 0.0199
$ws.Range("N6").Value = $newUpdate

# Row 7 - Curve (CRV)
$ws.Range("I7").Value = 2.87
$ws.Range("L7").Value = -0.0142
$ws.Range("N7").Value = $newUpdate

# Row 8 - Cover Protocol
$ws.Range("D8").Value = 665.0
$ws.Range("I8").Value = 487.37
$ws.Range("L8").Value = 0.006
$ws.Range("N8").Value = $newUpdate

# Row 9 - Ren (REN)
$ws.Range("I9").Value = 0.9268
$ws.Range("L9").Value = 0.0282
$ws.Range("N9").Value = $newUpdate

# Row 10 - Crypto.com Coin (CRO)
$ws.Range("I10").Value = 0.1941
$ws.Range("L10").Value = 0.0357
$ws.Range("N10").Value = $newUpdate
